# Add "극동미라주" (id 3420) as a new row, keeping the sheet's existing sort
# order by ID. It belongs right after "3418 왕십리금호어울림" (row 32) and
# before "3532 금호베스트빌" (the old row 33), so insert a fresh row 33 and
# push everything else down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

[void]$ws.Rows("33:33").Insert()
$ws.Range("A33").Value = 3420
$ws.Range("B33").Value = "극동미라주"

# Match the author's final cursor/selection state.
[void]$ws.Range("F29").Select()
